$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = '="60.230.69"'
$ws.Range('D2').Copy() | Out-Null
$ws.Range('D2').PasteSpecial(-4163) | Out-Null
$ws.Range('E2').Value = '  +4.49%  '

$ws.Range('D3').Formula = '="2.574.81"'
$ws.Range('D3').Copy() | Out-Null
$ws.Range('D3').PasteSpecial(-4163) | Out-Null
$ws.Range('E3').Value = '  +5.44%  '

$ws.Range('D4').Formula = '="0.999"'
$ws.Range('D4').Copy() | Out-Null
$ws.Range('D4').PasteSpecial(-4163) | Out-Null
$ws.Range('E4').Value = '  -0.30%  '

$ws.Range('D5').Formula = '="504.43"'
$ws.Range('D5').Copy() | Out-Null
$ws.Range('D5').PasteSpecial(-4163) | Out-Null
$ws.Range('E5').Value = '  +2.41%  '

$ws.Range('D6').Formula = '="154.52"'
$ws.Range('D6').Copy() | Out-Null
$ws.Range('D6').PasteSpecial(-4163) | Out-Null
$ws.Range('E6').Value = '  -1.57%  '

$ws.Range('D7').Formula = '="0.996"'
$ws.Range('D7').Copy() | Out-Null
$ws.Range('D7').PasteSpecial(-4163) | Out-Null
$ws.Range('E7').Value = '  +0.11%  '

$ws.Range('D8').Formula = '="0.574"'
$ws.Range('D8').Copy() | Out-Null
$ws.Range('D8').PasteSpecial(-4163) | Out-Null
$ws.Range('E8').Value = '  -6.65%  '

$ws.Range('D9').Formula = '="2.605.65"'
$ws.Range('D9').Copy() | Out-Null
$ws.Range('D9').PasteSpecial(-4163) | Out-Null
$ws.Range('E9').Value = '  +5.85%  '

$ws.Range('D10').Formula = '="6.46"'
$ws.Range('D10').Copy() | Out-Null
$ws.Range('D10').PasteSpecial(-4163) | Out-Null
$ws.Range('E10').Value = '  +2.58%  '

$ws.Range('D11').Formula = '="0.103"'
$ws.Range('D11').Copy() | Out-Null
$ws.Range('D11').PasteSpecial(-4163) | Out-Null
$ws.Range('E11').Value = '  +1.79%  '

$ws.Range('D12').Formula = '="0.340"'
$ws.Range('D12').Copy() | Out-Null
$ws.Range('D12').PasteSpecial(-4163) | Out-Null
$ws.Range('E12').Value = '  +1.28%  '

$ws.Range('E13').Value = '  +0.98%  '

$ws.Range('D14').Formula = '="3.038.54"'
$ws.Range('D14').Copy() | Out-Null
$ws.Range('D14').PasteSpecial(-4163) | Out-Null
$ws.Range('E14').Value = '  +5.97%  '

$ws.Range('D15').Formula = '="60.256.45"'
$ws.Range('D15').Copy() | Out-Null
$ws.Range('D15').PasteSpecial(-4163) | Out-Null
$ws.Range('E15').Value = '  +4.69%  '

$ws.Range('D16').Formula = '="21.57"'
$ws.Range('D16').Copy() | Out-Null
$ws.Range('D16').PasteSpecial(-4163) | Out-Null
$ws.Range('E16').Value = '  +3.27%  '

$ws.Range('E17').Value = '  +2.76%  '

$ws.Range('D18').Formula = '="2.604.41"'
$ws.Range('D18').Copy() | Out-Null
$ws.Range('D18').PasteSpecial(-4163) | Out-Null
$ws.Range('E18').Value = '  +5.61%  '

$ws.Range('D19').Formula = '="4.77"'
$ws.Range('D19').Copy() | Out-Null
$ws.Range('D19').PasteSpecial(-4163) | Out-Null
$ws.Range('E19').Value = '  +2.07%  '

$ws.Range('D20').Formula = '="339.97"'
$ws.Range('D20').Copy() | Out-Null
$ws.Range('D20').PasteSpecial(-4163) | Out-Null
$ws.Range('E20').Value = '  +3.87%  '

$ws.Range('D21').Formula = '="10.36"'
$ws.Range('D21').Copy() | Out-Null
$ws.Range('D21').PasteSpecial(-4163) | Out-Null
$ws.Range('E21').Value = '  +2.61%  '

$ws.Range('D22').Formula = '="6.07"'
$ws.Range('D22').Copy() | Out-Null
$ws.Range('D22').PasteSpecial(-4163) | Out-Null
$ws.Range('E22').Value = '  +1.52%  '

$ws.Range('E23').Value = '  +0.04%  '

$ws.Range('E24').Value = '  +2.84%  '

$ws.Range('D25').Formula = '="0.420"'
$ws.Range('D25').Copy() | Out-Null
$ws.Range('D25').PasteSpecial(-4163) | Out-Null
$ws.Range('E25').Value = '  +3.76%  '

$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Formula = '="2.704.98"'
$ws.Range('D26').Copy() | Out-Null
$ws.Range('D26').PasteSpecial(-4163) | Out-Null
$ws.Range('E26').Value = '  +5.81%  '

$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').Formula = '="0.165"'
$ws.Range('D27').Copy() | Out-Null
$ws.Range('D27').PasteSpecial(-4163) | Out-Null
$ws.Range('E27').Value = '  +2.38%  '

$ws.Range('D28').Formula = '="0.993"'
$ws.Range('D28').Copy() | Out-Null
$ws.Range('D28').PasteSpecial(-4163) | Out-Null
$ws.Range('E28').Value = '  -0.27%  '

$ws.Range('D29').Formula = '="0.0₃0852"'
$ws.Range('D29').Copy() | Out-Null
$ws.Range('D29').PasteSpecial(-4163) | Out-Null
$ws.Range('E29').Value = '  +5.03%  '

$ws.Range('D30').Formula = '="7.46"'
$ws.Range('D30').Copy() | Out-Null
$ws.Range('D30').PasteSpecial(-4163) | Out-Null
$ws.Range('E30').Value = '  +1.39%  '

$ws.Range('D31').Formula = '="0.998"'
$ws.Range('D31').Copy() | Out-Null
$ws.Range('D31').PasteSpecial(-4163) | Out-Null
$ws.Range('E31').Value = '  +0.01%  '

$ws.Range('D32').Formula = '="156.00"'
$ws.Range('D32').Copy() | Out-Null
$ws.Range('D32').PasteSpecial(-4163) | Out-Null
$ws.Range('E32').Value = '  +3.37%  '

$ws.Range('D33').Formula = '="19.25"'
$ws.Range('D33').Copy() | Out-Null
$ws.Range('D33').PasteSpecial(-4163) | Out-Null
$ws.Range('E33').Value = '  +2.11%  '

$ws.Range('E34').Value = '  +1.40%  '

$ws.Range('D35').Formula = '="5.70"'
$ws.Range('D35').Copy() | Out-Null
$ws.Range('D35').PasteSpecial(-4163) | Out-Null
$ws.Range('E35').Value = '  +6.38%  '

$ws.Range('D36').Formula = '="3.97"'
$ws.Range('D36').Copy() | Out-Null
$ws.Range('D36').PasteSpecial(-4163) | Out-Null
$ws.Range('E36').Value = '  +4.60%  '

$ws.Range('E37').Value = '  +4.51%  '

$ws.Range('D38').Formula = '="0.849"'
$ws.Range('D38').Copy() | Out-Null
$ws.Range('D38').PasteSpecial(-4163) | Out-Null
$ws.Range('E38').Value = '  +25.32%  '

$ws.Range('D39').Formula = '="3.76"'
$ws.Range('D39').Copy() | Out-Null
$ws.Range('D39').PasteSpecial(-4163) | Out-Null
$ws.Range('E39').Value = '  +5.56%  '

$ws.Range('D40').Formula = '="1.47"'
$ws.Range('D40').Copy() | Out-Null
$ws.Range('D40').PasteSpecial(-4163) | Out-Null
$ws.Range('E40').Value = '  +5.18%  '

$ws.Range('D41').Formula = '="0.840"'
$ws.Range('D41').Copy() | Out-Null
$ws.Range('D41').PasteSpecial(-4163) | Out-Null
$ws.Range('E41').Value = '  +0.32%  '

$ws.Range('D42').Formula = '="297.53"'
$ws.Range('D42').Copy() | Out-Null
$ws.Range('D42').PasteSpecial(-4163) | Out-Null
$ws.Range('E42').Value = '  +6.83%  '

$ws.Range('D43').Formula = '="35.48"'
$ws.Range('D43').Copy() | Out-Null
$ws.Range('D43').PasteSpecial(-4163) | Out-Null
$ws.Range('E43').Value = '  +3.41%  '

$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').Formula = '="0.0569"'
$ws.Range('D44').Copy() | Out-Null
$ws.Range('D44').PasteSpecial(-4163) | Out-Null
$ws.Range('E44').Value = '  +5.84%  '

$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').Formula = '="0.620"'
$ws.Range('D45').Copy() | Out-Null
$ws.Range('D45').PasteSpecial(-4163) | Out-Null
$ws.Range('E45').Value = '  +2.75%  '

$ws.Range('D46').Formula = '="0.0997"'
$ws.Range('D46').Copy() | Out-Null
$ws.Range('D46').PasteSpecial(-4163) | Out-Null
$ws.Range('E46').Value = '  -1.32%  '

$ws.Range('D47').Formula = '="0.996"'
$ws.Range('D47').Copy() | Out-Null
$ws.Range('D47').PasteSpecial(-4163) | Out-Null
$ws.Range('E47').Value = '  +0.22%  '

$ws.Range('D48').Formula = '="19.73"'
$ws.Range('D48').Copy() | Out-Null
$ws.Range('D48').PasteSpecial(-4163) | Out-Null
$ws.Range('E48').Value = '  +9.28%  '

$ws.Range('D49').Formula = '="4.90"'
$ws.Range('D49').Copy() | Out-Null
$ws.Range('D49').PasteSpecial(-4163) | Out-Null
$ws.Range('E49').Value = '  +4.33%  '

$ws.Range('D50').Formula = '="2.042.27"'
$ws.Range('D50').Copy() | Out-Null
$ws.Range('D50').PasteSpecial(-4163) | Out-Null
$ws.Range('E50').Value = '  +7.77%  '

$ws.Range('E51').Value = '  +1.10%  '
